# Applies the "automatic update of files" edit to the KLIPPAN overview sheet:
#  1) For every data row, bump the "Förändrad" date (column C) from 45184 to 45186.
#  2) For every HYPERLINK() formula that only has a URL argument (columns S, T, V,
#     W, X, Y), add the case's "Beteckning" (column A) as the friendly display
#     text, i.e. HYPERLINK("url") -> HYPERLINK("url", "Beteckning").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

$headerRow = 1
$firstDataRow = $headerRow + 1
$newChangedValue = 45186

$linkCols = @(19, 20, 22, 23, 24, 25)   # S, T, V, W, X, Y

for ($r = $firstDataRow; $r -le $lastRow; $r++) {

    # ---- 1) "Förändrad" date column (C = column 3) — every data row ----
    $ws.Cells.Item($r, 3).Value = $newChangedValue

    # ---- 2) Add friendly text to single-argument HYPERLINK formulas ----
    $beteckningCell = $ws.Cells.Item($r, 1)
    if ($beteckningCell.Value2 -eq $null) {
        continue
    }
    $beteckning = $beteckningCell.Value2

    foreach ($col in $linkCols) {
        $cell = $ws.Cells.Item($r, $col)
        if ($cell.HasFormula) {
            $f = $cell.Formula
            if ($f.StartsWith('=HYPERLINK(') -and -not $f.Contains('",')) {
                $updated = $f.Substring(0, $f.Length - 1) + ', "' + $beteckning + '")'
                $cell.Formula = $updated
            }
        }
    }
}

Write-Host "Update complete."
